$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.449.69'
$ws.Range("E2").Value = '  -1.48%  '
$ws.Range("D3").Value = '2.158.84'
$ws.Range("E3").Value = '  -3.28%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.39'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.605'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.49%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '71.64'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.87%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.574'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.59'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0902'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.84%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.01'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.11%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0999'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.66'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.43%  '
$ws.Range("D15").Value = '2.483.63'
$ws.Range("E15").Value = '  -3.25%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.09'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.59%  '
$ws.Range("D17").Value = '2.153.71'
$ws.Range("E17").Value = '  -3.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.776'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -7.45%  '
$ws.Range("D19").Value = '41.302.52'
$ws.Range("E19").Value = '  -1.52%  '
$ws.Range("E20").Value = '  -3.43%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '69.60'
$ws.Range("D21").Style = "Normal"
$ws.Range("E22").Value = '  -7.59%  '
$ws.Range("E23").Value = '  -12.96%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '226.69'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.02'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.62%  '
$ws.Range("E26").Value = '  +0.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.62'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -6.99%  '
$ws.Range("E28").Value = '  -8.93%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.18'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -4.63%  '
$ws.Range("E30").Value = '  -1.47%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '169.24'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.71'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.21%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '33.03'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +9.24%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0766'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.14'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -9.13%  '
$ws.Range("E36").Value = '  -4.04%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.28'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.63%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.104'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.45%  '
$ws.Range("E39").Value = '  -0.93%  '
$ws.Range("B40").Value = 'LidoDAOToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.08'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.41%  '
$ws.Range("B41").Value = 'Celestia'
$ws.Range("C41").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.86'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -12.01%  '
$ws.Range("E42").Value = '  -6.56%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '58.72'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -10.40%  '
$ws.Range("E44").Value = '  -5.44%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.36'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0957'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.26%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '95.78'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -8.35%  '
$ws.Range("E48").Value = '  -4.09%  '
$ws.Range("E49").Value = '  -5.76%  '
$ws.Range("E50").Value = '  -8.09%  '
$ws.Range("E51").Value = '  -2.41%  '
